# "Generate Report for handoff"
#
# The handoff transform for this file failed, so the generated report
# needs updating across the Overview sheet and the per-locale (zh-cn /
# de-de) sheets:
#   - Status text: "Ready for handoff" -> "Handoff transform failed"
#   - The "Latest Handoff File" cell (C2) + its hyperlink are cleared
#     (the transform never produced a handoff file).
#   - "Latest Handoff Datetime" (D2) resets to the zero/epoch datetime.
#   - "Handoff Reason" (H2) changes from "Include" to "Ignored".

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Helper: remove exactly the hyperlink anchored at $targetAddr (e.g. "$C$2")
# from a worksheet, leaving every other hyperlink untouched. Deleting while
# iterating the live Hyperlinks collection can skip an entry on a single
# pass, so keep retrying until the target address is actually gone.
function Remove-HyperlinkAt($ws, $targetAddr) {
    $guard = 0
    $stillThere = $true
    while ($stillThere -and $guard -lt 10) {
        foreach ($h in $ws.Hyperlinks) {
            if ($h.Range().Address() -eq $targetAddr) {
                $h.Delete()
                break
            }
        }
        $stillThere = $false
        foreach ($h in $ws.Hyperlinks) {
            if ($h.Range().Address() -eq $targetAddr) {
                $stillThere = $true
            }
        }
        $guard = $guard + 1
    }
}

# --- Update the status text everywhere it is used, so the old shared
#     string ("Ready for handoff") is fully replaced by the new one. ---
$overview.Range("B2").Value = "Handoff transform failed"
$overview.Range("C2").Value = "Handoff transform failed"
$zhcn.Range("B2").Value = "Handoff transform failed"
$dede.Range("B2").Value = "Handoff transform failed"

# --- zh-cn sheet: the handoff never produced a file, so clear it and
#     its hyperlink, reset the handoff datetime, and ignore it. ---
Remove-HyperlinkAt $zhcn '$C$2'
$zhcn.Range("C2").Clear()
$zhcn.Range("D2").Value = "0001-01-01 00:00:00"
$zhcn.Range("H2").Value = "Ignored"

# --- de-de sheet: same treatment. ---
Remove-HyperlinkAt $dede '$C$2'
$dede.Range("C2").Clear()
$dede.Range("D2").Value = "0001-01-01 00:00:00"
$dede.Range("H2").Value = "Ignored"
